$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2:AY2").ClearContents()
$ws.Range("A2").Value = 112038473
$ws.Range("B2").Value = 89686
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "Ovaliderad"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "NT"
$ws.Range("E2").Value = 658
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "Rosenticka"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "Rhodofomes roseus"
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "(Alb. & Schwein.) Kotl. & Pouzar"
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "4"
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "fruktkroppar"
$ws.Range("P2").NumberFormat = "@"
$ws.Range("P2").Value = "Lövnäs (Lövnäs), Ås lm"
$ws.Range("Q2").Value = 516057
$ws.Range("R2").Value = 7184320
$ws.Range("S2").Value = 10
$ws.Range("T2").NumberFormat = "@"
$ws.Range("T2").Value = "Västerbotten"
$ws.Range("U2").NumberFormat = "@"
$ws.Range("U2").Value = "Dorotea"
$ws.Range("V2").NumberFormat = "@"
$ws.Range("V2").Value = "Åsele lappmark"
$ws.Range("W2").NumberFormat = "@"
$ws.Range("W2").Value = "Dorotea"
$ws.Range("Y2").NumberFormat = "@"
$ws.Range("Y2").Value = "2023-09-11"
$ws.Range("Z2").NumberFormat = "@"
$ws.Range("Z2").Value = "13:34"
$ws.Range("AA2").NumberFormat = "@"
$ws.Range("AA2").Value = "2023-09-11"
$ws.Range("AB2").NumberFormat = "@"
$ws.Range("AB2").Value = "13:34"
$ws.Range("AD2").Value = $false
$ws.Range("AE2").Value = $false
$ws.Range("AG2").Value = $false
$ws.Range("AH2").NumberFormat = "@"
$ws.Range("AH2").Value = "Blåbärsgranskog"
$ws.Range("AJ2").NumberFormat = "@"
$ws.Range("AJ2").Value = "gran"
$ws.Range("AK2").NumberFormat = "@"
$ws.Range("AK2").Value = "Picea abies"
$ws.Range("AM2").NumberFormat = "@"
$ws.Range("AM2").Value = "Liggande död trädstam, markontakt"
$ws.Range("AO2").NumberFormat = "@"
$ws.Range("AO2").Value = "Horizontal, dead with ground contact # Picea abies"
$ws.Range("AW2").NumberFormat = "@"
$ws.Range("AW2").Value = "Eva Mårtensson"
$ws.Range("AX2").NumberFormat = "@"
$ws.Range("AX2").Value = "Eva Mårtensson"

# Row 3
$ws.Range("A3:AY3").ClearContents()
$ws.Range("A3").Value = 112037386
$ws.Range("B3").Value = 89423
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "Ovaliderad"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 5432
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "Granticka"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "Porodaedalea chrysoloma"
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("P3").NumberFormat = "@"
$ws.Range("P3").Value = "Lövnäs (Lövnäs), Ås lm"
$ws.Range("Q3").Value = 516032
$ws.Range("R3").Value = 7184227
$ws.Range("S3").Value = 10
$ws.Range("T3").NumberFormat = "@"
$ws.Range("T3").Value = "Västerbotten"
$ws.Range("U3").NumberFormat = "@"
$ws.Range("U3").Value = "Dorotea"
$ws.Range("V3").NumberFormat = "@"
$ws.Range("V3").Value = "Åsele lappmark"
$ws.Range("W3").NumberFormat = "@"
$ws.Range("W3").Value = "Dorotea"
$ws.Range("Y3").NumberFormat = "@"
$ws.Range("Y3").Value = "2023-09-11"
$ws.Range("Z3").NumberFormat = "@"
$ws.Range("Z3").Value = "11:52"
$ws.Range("AA3").NumberFormat = "@"
$ws.Range("AA3").Value = "2023-09-11"
$ws.Range("AB3").NumberFormat = "@"
$ws.Range("AB3").Value = "11:52"
$ws.Range("AD3").Value = $false
$ws.Range("AE3").Value = $false
$ws.Range("AG3").Value = $false
$ws.Range("AH3").NumberFormat = "@"
$ws.Range("AH3").Value = "Blåbärsbarrskog"
$ws.Range("AJ3").NumberFormat = "@"
$ws.Range("AJ3").Value = "gran"
$ws.Range("AK3").NumberFormat = "@"
$ws.Range("AK3").Value = "Picea abies"
$ws.Range("AM3").NumberFormat = "@"
$ws.Range("AM3").Value = "Liggande död trädstam, markontakt"
$ws.Range("AO3").NumberFormat = "@"
$ws.Range("AO3").Value = "Horizontal, dead with ground contact # Picea abies"
$ws.Range("AW3").NumberFormat = "@"
$ws.Range("AW3").Value = "Eva Mårtensson"
$ws.Range("AX3").NumberFormat = "@"
$ws.Range("AX3").Value = "Eva Mårtensson"

# Row 4
$ws.Range("A4:AY4").ClearContents()
$ws.Range("A4").Value = 112038134
$ws.Range("B4").Value = 89405
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "Ovaliderad"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 1202
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "Ullticka"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("P4").NumberFormat = "@"
$ws.Range("P4").Value = "Lövnäs (Lövnäs), Ås lm"
$ws.Range("Q4").Value = 515925
$ws.Range("R4").Value = 7184319
$ws.Range("S4").Value = 10
$ws.Range("T4").NumberFormat = "@"
$ws.Range("T4").Value = "Västerbotten"
$ws.Range("U4").NumberFormat = "@"
$ws.Range("U4").Value = "Dorotea"
$ws.Range("V4").NumberFormat = "@"
$ws.Range("V4").Value = "Åsele lappmark"
$ws.Range("W4").NumberFormat = "@"
$ws.Range("W4").Value = "Dorotea"
$ws.Range("Y4").NumberFormat = "@"
$ws.Range("Y4").Value = "2023-09-11"
$ws.Range("Z4").NumberFormat = "@"
$ws.Range("Z4").Value = "13:27"
$ws.Range("AA4").NumberFormat = "@"
$ws.Range("AA4").Value = "2023-09-11"
$ws.Range("AB4").NumberFormat = "@"
$ws.Range("AB4").Value = "13:27"
$ws.Range("AD4").Value = $false
$ws.Range("AE4").Value = $false
$ws.Range("AG4").Value = $false
$ws.Range("AH4").NumberFormat = "@"
$ws.Range("AH4").Value = "Blåbärsgranskog"
$ws.Range("AJ4").NumberFormat = "@"
$ws.Range("AJ4").Value = "gran"
$ws.Range("AK4").NumberFormat = "@"
$ws.Range("AK4").Value = "Picea abies"
$ws.Range("AM4").NumberFormat = "@"
$ws.Range("AM4").Value = "Liggande död trädstam, markontakt"
$ws.Range("AO4").NumberFormat = "@"
$ws.Range("AO4").Value = "Horizontal, dead with ground contact # Picea abies"
$ws.Range("AW4").NumberFormat = "@"
$ws.Range("AW4").Value = "Eva Mårtensson"
$ws.Range("AX4").NumberFormat = "@"
$ws.Range("AX4").Value = "Eva Mårtensson"

# Row 5
$ws.Range("A5:AY5").ClearContents()
$ws.Range("A5").Value = 112037208
$ws.Range("B5").Value = 77515
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "Ovaliderad"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 6425
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "Garnlav"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "Alectoria sarmentosa"
$ws.Range("H5").NumberFormat = "@"
$ws.Range("H5").Value = "(Ach.) Ach."
$ws.Range("P5").NumberFormat = "@"
$ws.Range("P5").Value = "Lövnäs (Lövnäs), Ås lm"
$ws.Range("Q5").Value = 516097
$ws.Range("R5").Value = 7184259
$ws.Range("S5").Value = 10
$ws.Range("T5").NumberFormat = "@"
$ws.Range("T5").Value = "Västerbotten"
$ws.Range("U5").NumberFormat = "@"
$ws.Range("U5").Value = "Dorotea"
$ws.Range("V5").NumberFormat = "@"
$ws.Range("V5").Value = "Åsele lappmark"
$ws.Range("W5").NumberFormat = "@"
$ws.Range("W5").Value = "Dorotea"
$ws.Range("Y5").NumberFormat = "@"
$ws.Range("Y5").Value = "2023-09-11"
$ws.Range("Z5").NumberFormat = "@"
$ws.Range("Z5").Value = "11:44"
$ws.Range("AA5").NumberFormat = "@"
$ws.Range("AA5").Value = "2023-09-11"
$ws.Range("AB5").NumberFormat = "@"
$ws.Range("AB5").Value = "11:44"
$ws.Range("AD5").Value = $false
$ws.Range("AE5").Value = $false
$ws.Range("AG5").Value = $false
$ws.Range("AH5").NumberFormat = "@"
$ws.Range("AH5").Value = "Blåbärsbarrskog"
$ws.Range("AM5").NumberFormat = "@"
$ws.Range("AM5").Value = "Stående död trädstam/högstubbe"
$ws.Range("AO5").NumberFormat = "@"
$ws.Range("AO5").Value = "Standing dead tree/snags"
$ws.Range("AW5").NumberFormat = "@"
$ws.Range("AW5").Value = "Eva Mårtensson"
$ws.Range("AX5").NumberFormat = "@"
$ws.Range("AX5").Value = "Eva Mårtensson"

# Row 6
$ws.Range("A6:AY6").ClearContents()
$ws.Range("A6").Value = 112038436
$ws.Range("B6").Value = 89401
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "Ovaliderad"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "NT"
$ws.Range("E6").Value = 1108
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = "Harticka"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "Pelloporus leporinus"
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = "(Fr.) Krieglst."
$ws.Range("P6").NumberFormat = "@"
$ws.Range("P6").Value = "Lövnäs (Lövnäs), Ås lm"
$ws.Range("Q6").Value = 515951
$ws.Range("R6").Value = 7184320
$ws.Range("S6").Value = 10
$ws.Range("T6").NumberFormat = "@"
$ws.Range("T6").Value = "Västerbotten"
$ws.Range("U6").NumberFormat = "@"
$ws.Range("U6").Value = "Dorotea"
$ws.Range("V6").NumberFormat = "@"
$ws.Range("V6").Value = "Åsele lappmark"
$ws.Range("W6").NumberFormat = "@"
$ws.Range("W6").Value = "Dorotea"
$ws.Range("Y6").NumberFormat = "@"
$ws.Range("Y6").Value = "2023-09-11"
$ws.Range("Z6").NumberFormat = "@"
$ws.Range("Z6").Value = "13:28"
$ws.Range("AA6").NumberFormat = "@"
$ws.Range("AA6").Value = "2023-09-11"
$ws.Range("AB6").NumberFormat = "@"
$ws.Range("AB6").Value = "13:28"
$ws.Range("AD6").Value = $false
$ws.Range("AE6").Value = $false
$ws.Range("AG6").Value = $false
$ws.Range("AH6").NumberFormat = "@"
$ws.Range("AH6").Value = "Blåbärsgranskog"
$ws.Range("AJ6").NumberFormat = "@"
$ws.Range("AJ6").Value = "gran"
$ws.Range("AK6").NumberFormat = "@"
$ws.Range("AK6").Value = "Picea abies"
$ws.Range("AM6").NumberFormat = "@"
$ws.Range("AM6").Value = "Stående död trädstam/högstubbe"
$ws.Range("AO6").NumberFormat = "@"
$ws.Range("AO6").Value = "Standing dead tree/snags # Picea abies"
$ws.Range("AW6").NumberFormat = "@"
$ws.Range("AW6").Value = "Eva Mårtensson"
$ws.Range("AX6").NumberFormat = "@"
$ws.Range("AX6").Value = "Eva Mårtensson"

# Row 7
$ws.Range("A7:AY7").ClearContents()
$ws.Range("A7").Value = 112035549
$ws.Range("B7").Value = 77515
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "Ovaliderad"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "NT"
$ws.Range("E7").Value = 6425
$ws.Range("F7").NumberFormat = "@"
$ws.Range("F7").Value = "Garnlav"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "Alectoria sarmentosa"
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = "(Ach.) Ach."
$ws.Range("P7").NumberFormat = "@"
$ws.Range("P7").Value = "Lövnäs (Lövnäs), Ås lm"
$ws.Range("Q7").Value = 515977
$ws.Range("R7").Value = 7184567
$ws.Range("S7").Value = 10
$ws.Range("T7").NumberFormat = "@"
$ws.Range("T7").Value = "Västerbotten"
$ws.Range("U7").NumberFormat = "@"
$ws.Range("U7").Value = "Dorotea"
$ws.Range("V7").NumberFormat = "@"
$ws.Range("V7").Value = "Åsele lappmark"
$ws.Range("W7").NumberFormat = "@"
$ws.Range("W7").Value = "Dorotea"
$ws.Range("Y7").NumberFormat = "@"
$ws.Range("Y7").Value = "2023-09-11"
$ws.Range("Z7").NumberFormat = "@"
$ws.Range("Z7").Value = "10:51"
$ws.Range("AA7").NumberFormat = "@"
$ws.Range("AA7").Value = "2023-09-11"
$ws.Range("AB7").NumberFormat = "@"
$ws.Range("AB7").Value = "10:51"
$ws.Range("AD7").Value = $false
$ws.Range("AE7").Value = $false
$ws.Range("AG7").Value = $false
$ws.Range("AH7").NumberFormat = "@"
$ws.Range("AH7").Value = "Blåbärsgranskog"
$ws.Range("AM7").NumberFormat = "@"
$ws.Range("AM7").Value = "Gren på levande träd"
$ws.Range("AO7").NumberFormat = "@"
$ws.Range("AO7").Value = "Branch on living tree"
$ws.Range("AW7").NumberFormat = "@"
$ws.Range("AW7").Value = "Eva Mårtensson"
$ws.Range("AX7").NumberFormat = "@"
$ws.Range("AX7").Value = "Eva Mårtensson"

# Row 8
$ws.Range("A8:AY8").ClearContents()
$ws.Range("A8").Value = 112038529
$ws.Range("B8").Value = 77515
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "Ovaliderad"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "NT"
$ws.Range("E8").Value = 6425
$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = "Garnlav"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "Alectoria sarmentosa"
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = "(Ach.) Ach."
$ws.Range("P8").NumberFormat = "@"
$ws.Range("P8").Value = "Lövnäs (Lövnäs), Ås lm"
$ws.Range("Q8").Value = 515872
$ws.Range("R8").Value = 7184628
$ws.Range("S8").Value = 10
$ws.Range("T8").NumberFormat = "@"
$ws.Range("T8").Value = "Västerbotten"
$ws.Range("U8").NumberFormat = "@"
$ws.Range("U8").Value = "Dorotea"
$ws.Range("V8").NumberFormat = "@"
$ws.Range("V8").Value = "Åsele lappmark"
$ws.Range("W8").NumberFormat = "@"
$ws.Range("W8").Value = "Dorotea"
$ws.Range("Y8").NumberFormat = "@"
$ws.Range("Y8").Value = "2023-09-11"
$ws.Range("Z8").NumberFormat = "@"
$ws.Range("Z8").Value = "14:07"
$ws.Range("AA8").NumberFormat = "@"
$ws.Range("AA8").Value = "2023-09-11"
$ws.Range("AB8").NumberFormat = "@"
$ws.Range("AB8").Value = "14:07"
$ws.Range("AD8").Value = $false
$ws.Range("AE8").Value = $false
$ws.Range("AG8").Value = $false
$ws.Range("AH8").NumberFormat = "@"
$ws.Range("AH8").Value = "Gransumpskog"
$ws.Range("AW8").NumberFormat = "@"
$ws.Range("AW8").Value = "Eva Mårtensson"
$ws.Range("AX8").NumberFormat = "@"
$ws.Range("AX8").Value = "Eva Mårtensson"

# Row 9
$ws.Range("A9:AY9").ClearContents()
$ws.Range("A9").Value = 112035020
$ws.Range("B9").Value = 89401
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "Ovaliderad"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "NT"
$ws.Range("E9").Value = 1108
$ws.Range("F9").NumberFormat = "@"
$ws.Range("F9").Value = "Harticka"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "Pelloporus leporinus"
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "(Fr.) Krieglst."
$ws.Range("P9").NumberFormat = "@"
$ws.Range("P9").Value = "Lövnäs (Lövnäs), Ås lm"
$ws.Range("Q9").Value = 515923
$ws.Range("R9").Value = 7184659
$ws.Range("S9").Value = 50
$ws.Range("T9").NumberFormat = "@"
$ws.Range("T9").Value = "Västerbotten"
$ws.Range("U9").NumberFormat = "@"
$ws.Range("U9").Value = "Dorotea"
$ws.Range("V9").NumberFormat = "@"
$ws.Range("V9").Value = "Åsele lappmark"
$ws.Range("W9").NumberFormat = "@"
$ws.Range("W9").Value = "Dorotea"
$ws.Range("Y9").NumberFormat = "@"
$ws.Range("Y9").Value = "2023-09-11"
$ws.Range("Z9").NumberFormat = "@"
$ws.Range("Z9").Value = "10:24"
$ws.Range("AA9").NumberFormat = "@"
$ws.Range("AA9").Value = "2023-09-11"
$ws.Range("AB9").NumberFormat = "@"
$ws.Range("AB9").Value = "10:24"
$ws.Range("AD9").Value = $false
$ws.Range("AE9").Value = $false
$ws.Range("AG9").Value = $false
$ws.Range("AH9").NumberFormat = "@"
$ws.Range("AH9").Value = "Blåbärsbarrskog"
$ws.Range("AM9").NumberFormat = "@"
$ws.Range("AM9").Value = "Stubbe"
$ws.Range("AO9").NumberFormat = "@"
$ws.Range("AO9").Value = "Stump"
$ws.Range("AW9").NumberFormat = "@"
$ws.Range("AW9").Value = "Eva Mårtensson"
$ws.Range("AX9").NumberFormat = "@"
$ws.Range("AX9").Value = "Eva Mårtensson"

# Row 10
$ws.Range("A10:AY10").ClearContents()
$ws.Range("A10").Value = 112037635
$ws.Range("B10").Value = 89401
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "Ovaliderad"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "NT"
$ws.Range("E10").Value = 1108
$ws.Range("F10").NumberFormat = "@"
$ws.Range("F10").Value = "Harticka"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "Pelloporus leporinus"
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H10").Value = "(Fr.) Krieglst."
$ws.Range("P10").NumberFormat = "@"
$ws.Range("P10").Value = "Lövnäs (Lövnäs), Ås lm"
$ws.Range("Q10").Value = 515886
$ws.Range("R10").Value = 7184226
$ws.Range("S10").Value = 10
$ws.Range("T10").NumberFormat = "@"
$ws.Range("T10").Value = "Västerbotten"
$ws.Range("U10").NumberFormat = "@"
$ws.Range("U10").Value = "Dorotea"
$ws.Range("V10").NumberFormat = "@"
$ws.Range("V10").Value = "Åsele lappmark"
$ws.Range("W10").NumberFormat = "@"
$ws.Range("W10").Value = "Dorotea"
$ws.Range("Y10").NumberFormat = "@"
$ws.Range("Y10").Value = "2023-09-11"
$ws.Range("Z10").NumberFormat = "@"
$ws.Range("Z10").Value = "12:06"
$ws.Range("AA10").NumberFormat = "@"
$ws.Range("AA10").Value = "2023-09-11"
$ws.Range("AB10").NumberFormat = "@"
$ws.Range("AB10").Value = "12:06"
$ws.Range("AD10").Value = $false
$ws.Range("AE10").Value = $false
$ws.Range("AG10").Value = $false
$ws.Range("AH10").NumberFormat = "@"
$ws.Range("AH10").Value = "Blåbärsgranskog"
$ws.Range("AJ10").NumberFormat = "@"
$ws.Range("AJ10").Value = "gran"
$ws.Range("AK10").NumberFormat = "@"
$ws.Range("AK10").Value = "Picea abies"
$ws.Range("AM10").NumberFormat = "@"
$ws.Range("AM10").Value = "Liggande död trädstam, markontakt"
$ws.Range("AO10").NumberFormat = "@"
$ws.Range("AO10").Value = "Horizontal, dead with ground contact # Picea abies"
$ws.Range("AW10").NumberFormat = "@"
$ws.Range("AW10").Value = "Eva Mårtensson"
$ws.Range("AX10").NumberFormat = "@"
$ws.Range("AX10").Value = "Eva Mårtensson"

# Row 11
$ws.Range("A11:AY11").ClearContents()
$ws.Range("A11").Value = 112037684
$ws.Range("B11").Value = 77515
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "Ovaliderad"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "NT"
$ws.Range("E11").Value = 6425
$ws.Range("F11").NumberFormat = "@"
$ws.Range("F11").Value = "Garnlav"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "Alectoria sarmentosa"
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = "(Ach.) Ach."
$ws.Range("P11").NumberFormat = "@"
$ws.Range("P11").Value = "Lövnäs (Lövnäs), Ås lm"
$ws.Range("Q11").Value = 515886
$ws.Range("R11").Value = 7184226
$ws.Range("S11").Value = 10
$ws.Range("T11").NumberFormat = "@"
$ws.Range("T11").Value = "Västerbotten"
$ws.Range("U11").NumberFormat = "@"
$ws.Range("U11").Value = "Dorotea"
$ws.Range("V11").NumberFormat = "@"
$ws.Range("V11").Value = "Åsele lappmark"
$ws.Range("W11").NumberFormat = "@"
$ws.Range("W11").Value = "Dorotea"
$ws.Range("Y11").NumberFormat = "@"
$ws.Range("Y11").Value = "2023-09-11"
$ws.Range("Z11").NumberFormat = "@"
$ws.Range("Z11").Value = "12:08"
$ws.Range("AA11").NumberFormat = "@"
$ws.Range("AA11").Value = "2023-09-11"
$ws.Range("AB11").NumberFormat = "@"
$ws.Range("AB11").Value = "12:08"
$ws.Range("AD11").Value = $false
$ws.Range("AE11").Value = $false
$ws.Range("AG11").Value = $false
$ws.Range("AH11").NumberFormat = "@"
$ws.Range("AH11").Value = "Blåbärsbarrskog"
$ws.Range("AJ11").NumberFormat = "@"
$ws.Range("AJ11").Value = "gran"
$ws.Range("AK11").NumberFormat = "@"
$ws.Range("AK11").Value = "Picea abies"
$ws.Range("AM11").NumberFormat = "@"
$ws.Range("AM11").Value = "Stående död trädstam/högstubbe"
$ws.Range("AO11").NumberFormat = "@"
$ws.Range("AO11").Value = "Standing dead tree/snags # Picea abies"
$ws.Range("AW11").NumberFormat = "@"
$ws.Range("AW11").Value = "Eva Mårtensson"
$ws.Range("AX11").NumberFormat = "@"
$ws.Range("AX11").Value = "Eva Mårtensson"

# Row 12
$ws.Range("A12:AY12").ClearContents()
$ws.Range("A12").Value = 112035981
$ws.Range("B12").Value = 90687
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "Ovaliderad"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "LC"
$ws.Range("E12").Value = 5964
$ws.Range("F12").NumberFormat = "@"
$ws.Range("F12").Value = "Fjällig taggsvamp s.str."
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "Sarcodon imbricatus s.str."
$ws.Range("H12").NumberFormat = "@"
$ws.Range("H12").Value = "(L.:Fr.) P.Karst."
$ws.Range("P12").NumberFormat = "@"
$ws.Range("P12").Value = "Lövnäs (Lövnäs), Ås lm"
$ws.Range("Q12").Value = 516149
$ws.Range("R12").Value = 7184413
$ws.Range("S12").Value = 5
$ws.Range("T12").NumberFormat = "@"
$ws.Range("T12").Value = "Västerbotten"
$ws.Range("U12").NumberFormat = "@"
$ws.Range("U12").Value = "Dorotea"
$ws.Range("V12").NumberFormat = "@"
$ws.Range("V12").Value = "Åsele lappmark"
$ws.Range("W12").NumberFormat = "@"
$ws.Range("W12").Value = "Dorotea"
$ws.Range("Y12").NumberFormat = "@"
$ws.Range("Y12").Value = "2023-09-11"
$ws.Range("Z12").NumberFormat = "@"
$ws.Range("Z12").Value = "11:29"
$ws.Range("AA12").NumberFormat = "@"
$ws.Range("AA12").Value = "2023-09-11"
$ws.Range("AB12").NumberFormat = "@"
$ws.Range("AB12").Value = "11:29"
$ws.Range("AD12").Value = $false
$ws.Range("AE12").Value = $false
$ws.Range("AG12").Value = $false
$ws.Range("AH12").NumberFormat = "@"
$ws.Range("AH12").Value = "Blåbärsbarrskog"
$ws.Range("AW12").NumberFormat = "@"
$ws.Range("AW12").Value = "Eva Mårtensson"
$ws.Range("AX12").NumberFormat = "@"
$ws.Range("AX12").Value = "Eva Mårtensson"

# Row 13
$ws.Range("A13:AY13").ClearContents()
$ws.Range("A13").Value = 112038082
$ws.Range("B13").Value = 90087
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "Ovaliderad"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "LC"
$ws.Range("E13").Value = 3298
$ws.Range("F13").NumberFormat = "@"
$ws.Range("F13").Value = "Trådticka"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "Climacocystis borealis"
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = "(Fr.) Kotl. & Pouzar"
$ws.Range("P13").NumberFormat = "@"
$ws.Range("P13").Value = "Lövnäs (Lövnäs), Ås lm"
$ws.Range("Q13").Value = 515925
$ws.Range("R13").Value = 7184319
$ws.Range("S13").Value = 10
$ws.Range("T13").NumberFormat = "@"
$ws.Range("T13").Value = "Västerbotten"
$ws.Range("U13").NumberFormat = "@"
$ws.Range("U13").Value = "Dorotea"
$ws.Range("V13").NumberFormat = "@"
$ws.Range("V13").Value = "Åsele lappmark"
$ws.Range("W13").NumberFormat = "@"
$ws.Range("W13").Value = "Dorotea"
$ws.Range("Y13").NumberFormat = "@"
$ws.Range("Y13").Value = "2023-09-11"
$ws.Range("Z13").NumberFormat = "@"
$ws.Range("Z13").Value = "13:22"
$ws.Range("AA13").NumberFormat = "@"
$ws.Range("AA13").Value = "2023-09-11"
$ws.Range("AB13").NumberFormat = "@"
$ws.Range("AB13").Value = "13:22"
$ws.Range("AD13").Value = $false
$ws.Range("AE13").Value = $false
$ws.Range("AG13").Value = $false
$ws.Range("AH13").NumberFormat = "@"
$ws.Range("AH13").Value = "Blåbärsgranskog"
$ws.Range("AJ13").NumberFormat = "@"
$ws.Range("AJ13").Value = "gran"
$ws.Range("AK13").NumberFormat = "@"
$ws.Range("AK13").Value = "Picea abies"
$ws.Range("AM13").NumberFormat = "@"
$ws.Range("AM13").Value = "Stående död trädstam/högstubbe"
$ws.Range("AO13").NumberFormat = "@"
$ws.Range("AO13").Value = "Standing dead tree/snags # Picea abies"
$ws.Range("AW13").NumberFormat = "@"
$ws.Range("AW13").Value = "Eva Mårtensson"
$ws.Range("AX13").NumberFormat = "@"
$ws.Range("AX13").Value = "Eva Mårtensson"

Write-Host "Edit complete"